$wb = $excel.ActiveWorkbook
$main = $wb.Worksheets.Item("Test Cases")

# --- Rename TestCase_A# -> IAM0## in the "Test Cases" summary sheet, column A ---
# Row r holds TestCase_A(r-1). The edit order below (13 then 20 then 14..19 then 21..26)
# matches how the shared-strings table ends up ordered in the saved file.
$rowOrder = @(2,3,4,5,6,7,8,9,10,11,12,13,20,14,15,16,17,18,19,21,22,23,24,25,26)
foreach ($row in $rowOrder) {
    $n = $row - 1
    $newName = "IAM{0:D3}" -f $n
    $main.Cells.Item($row, 1).Value = $newName
}

# --- Column header E1 changed from "Results" to "PASS" ---
$main.Range("E1").Value = "PASS"

# --- Fix inconsistent cell style on A22 (align with neighboring rows) ---
$main.Range("A22").Style = $main.Range("A21").Style

# --- Rename the TestCase_A* sheet tabs to IAM0## ---
$renames = @(
    @("TestCase_A5",  "IAM005"),
    @("TestCase_A6",  "IAM006"),
    @("TestCase_A7",  "IAM007"),
    @("TestCase_A8",  "IAM008"),
    @("TestCase_A9",  "IAM009"),
    @("TestCase_A10", "IAM010"),
    @("TestCase_A11", "IAM011"),
    @("TestCase_A12", "IAM012"),
    @("TestCase_A19", "IAM019")
)
foreach ($pair in $renames) {
    $ws = $wb.Worksheets.Item($pair[0])
    $ws.Name = $pair[1]
}

# --- Cosmetic: restore the active-cell selections that moved in the source sheets ---
$iam005 = $wb.Worksheets.Item("IAM005")
$iam005.Range("L7").Select()

$iam006 = $wb.Worksheets.Item("IAM006")
$iam006.Range("D33").Select()

$iam012 = $wb.Worksheets.Item("IAM012")
$iam012.Range("B38").Select()

$main.Select()
